# Auto-generated Excel COM-interop script applying the commit diff
# "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3500
$ws.Range("I106").Value = 3500
$ws.Range("K106").Value = 3500
$ws.Range("M106").Value = -2869

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3030.6296
$ws.Range("I137").Value = 2625.5334
$ws.Range("K137").Value = 7876.600199999999
$ws.Range("M137").Value = -5326.600199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1670614.2
$ws.Range("I32").Value = 1740197.5
$ws.Range("K32").Value = 1740197.5
$ws.Range("M32").Value = -1739910.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5373.5
$ws.Range("I45").Value = 2746.889
$ws.Range("K45").Value = 2746.889
$ws.Range("M45").Value = -2369.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5002.1523
$ws.Range("I61").Value = 2033.8235
$ws.Range("K61").Value = 2033.8235
$ws.Range("M61").Value = -1821.8235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2597.111
$ws.Range("I63").Value = 2369.75
$ws.Range("K63").Value = 2369.75
$ws.Range("M63").Value = -1683.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2597.111
$ws.Range("I66").Value = 2369.75
$ws.Range("K66").Value = 11848.75
$ws.Range("M66").Value = -8416.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1264.6875
$ws.Range("I102").Value = 1245.3572
$ws.Range("K102").Value = 1245.3572
$ws.Range("M102").Value = 376.6428000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 55557084
$ws.Range("I110").Value = 1705.5
$ws.Range("J110").Value = 83334776
$ws.Range("K110").Value = 1705.5
$ws.Range("L110").Value = 83334776
$ws.Range("M110").Value = 339.5
$ws.Range("N110").Value = -83338866

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 24432
$ws.Range("I122").Value = 30331.428
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 90994.284
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -88544.284
$ws.Range("N122").Value = -36900.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5002.1523
$ws.Range("I136").Value = 2033.8235
$ws.Range("K136").Value = 6101.470499999999
$ws.Range("M136").Value = -3551.470499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5433.367
$ws.Range("I134").Value = 2025.1154
$ws.Range("K134").Value = 6075.3462
$ws.Range("M134").Value = -3540.3462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 391.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6922.3076
$ws.Range("I99").Value = 6997.923
$ws.Range("J99").Value = 6846.6924
$ws.Range("K99").Value = 6997.923
$ws.Range("L99").Value = 6846.6924
$ws.Range("M99").Value = -5499.923
$ws.Range("N99").Value = -9842.6924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6922.3076
$ws.Range("I126").Value = 6997.923
$ws.Range("J126").Value = 6846.6924
$ws.Range("K126").Value = 20993.769
$ws.Range("L126").Value = 20540.0772
$ws.Range("M126").Value = -18523.769
$ws.Range("N126").Value = -25480.0772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9095263
$ws.Range("I132").Value = 1315.3077
$ws.Range("J132").Value = 22230966
$ws.Range("K132").Value = 3945.9231
$ws.Range("L132").Value = 66692898
$ws.Range("M132").Value = -1415.9231
$ws.Range("N132").Value = -66697958

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 336447870
$ws.Range("J4").Value = 336228450
$ws.Range("L4").Value = 1008685350
$ws.Range("N4").Value = -1008685574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2354789
$ws.Range("J5").Value = 2815.1
$ws.Range("L5").Value = 8445.299999999999
$ws.Range("N5").Value = -8669.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 67
$ws.Range("J38").Value = 63.6
$ws.Range("L38").Value = 190.8
$ws.Range("N38").Value = -884.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 12500260
$ws.Range("I107").Value = 315.75
$ws.Range("J107").Value = 16666908
$ws.Range("K107").Value = 947.25
$ws.Range("L107").Value = 50000724
$ws.Range("M107").Value = 972.75
$ws.Range("N107").Value = -50004564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2354789
$ws.Range("J135").Value = 2815.1
$ws.Range("L135").Value = 25335.9
$ws.Range("N135").Value = -30405.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 113148.945
$ws.Range("I137").Value = 112648.22
$ws.Range("J137").Value = 113649.664
$ws.Range("K137").Value = 337944.66
$ws.Range("L137").Value = 340948.992
$ws.Range("M137").Value = -332844.66
$ws.Range("N137").Value = -351148.992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3054
$ws.Range("I113").Value = 3011.4443
$ws.Range("K113").Value = 3011.4443
$ws.Range("M113").Value = -841.4443000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3574.5
$ws.Range("I126").Value = 3794.6
$ws.Range("J126").Value = 2474
$ws.Range("K126").Value = 11383.8
$ws.Range("L126").Value = 7422
$ws.Range("M126").Value = -8913.799999999999
$ws.Range("N126").Value = -12362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4896.722
$ws.Range("I132").Value = 1180.8
$ws.Range("J132").Value = 9541.625
$ws.Range("K132").Value = 3542.4
$ws.Range("L132").Value = 28624.875
$ws.Range("M132").Value = -1012.4
$ws.Range("N132").Value = -33684.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 15380
$ws.Range("I93").Value = 13500
$ws.Range("J93").Value = 16633.334
$ws.Range("K93").Value = 13500
$ws.Range("L93").Value = 16633.334
$ws.Range("M93").Value = -12252
$ws.Range("N93").Value = -19129.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10422915
$ws.Range("I132").Value = 23811880
$ws.Range("J132").Value = 9274.037
$ws.Range("K132").Value = 71435640
$ws.Range("L132").Value = 27822.111
$ws.Range("M132").Value = -71433110
$ws.Range("N132").Value = -32882.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8873.441000000001
$ws.Range("I136").Value = 2558.0833
$ws.Range("J136").Value = 12318.182
$ws.Range("K136").Value = 7674.249899999999
$ws.Range("L136").Value = 36954.546
$ws.Range("M136").Value = -5124.249899999999
$ws.Range("N136").Value = -42054.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 4811250
$ws.Range("I5").Value = 122500
$ws.Range("J5").Value = 9500000
$ws.Range("K5").Value = 122500
$ws.Range("L5").Value = 9500000
$ws.Range("M5").Value = -122388
$ws.Range("N5").Value = -9500224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 15005
$ws.Range("J69").Value = 15005
$ws.Range("L69").Value = 15005
$ws.Range("N69").Value = -16503

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 15005
$ws.Range("J72").Value = 15005
$ws.Range("L72").Value = 45015
$ws.Range("N72").Value = -52503

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 15391189
$ws.Range("I81").Value = 2314.2727
$ws.Range("K81").Value = 4628.5454
$ws.Range("M81").Value = -3567.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 15391189
$ws.Range("I84").Value = 2314.2727
$ws.Range("K84").Value = 23142.727
$ws.Range("M84").Value = -17838.727
